$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 15000
$ws.Range("I21").Value = 15000
$ws.Range("K21").Value = 15000
$ws.Range("M21").Value = -14532

$ws.Range("H23").Value = 15000
$ws.Range("I23").Value = 15000
$ws.Range("K23").Value = 15000
$ws.Range("M23").Value = -14766

$ws.Range("H80").Value = 877033.75
$ws.Range("I80").Value = 1422318.8
$ws.Range("K80").Value = 4266956.4
$ws.Range("M80").Value = -4265958.4

$ws.Range("H83").Value = 877033.75
$ws.Range("I83").Value = 1422318.8
$ws.Range("K83").Value = 12800869.2
$ws.Range("M83").Value = -12795877.2

$ws.Range("H86").Value = 5292438.5
$ws.Range("I86").Value = 2600.875
$ws.Range("K86").Value = 2600.875
$ws.Range("M86").Value = -1477.875

$ws.Range("H88").Value = 2109.389
$ws.Range("I88").Value = 1529.8334
$ws.Range("K88").Value = 1529.8334
$ws.Range("M88").Value = -1123.8334

$ws.Range("H89").Value = 5292438.5
$ws.Range("I89").Value = 2600.875
$ws.Range("K89").Value = 13004.375
$ws.Range("M89").Value = -7388.375

$ws.Range("H91").Value = 2109.389
$ws.Range("I91").Value = 1529.8334
$ws.Range("K91").Value = 1529.8334
$ws.Range("M91").Value = -125.8334

$ws.Range("H92").Value = 117324.47
$ws.Range("J92").Value = 331728.5
$ws.Range("L92").Value = 331728.5
$ws.Range("N92").Value = -334224.5

$ws.Range("H96").Value = 453.42856
$ws.Range("I96").Value = 299.66666
$ws.Range("J96").Value = 568.75
$ws.Range("K96").Value = 898.9999799999999
$ws.Range("L96").Value = 1706.25
$ws.Range("M96").Value = 474.0000200000001
$ws.Range("N96").Value = -4452.25

$ws.Range("H97").Value = 496.5
$ws.Range("J97").Value = 495.33334
$ws.Range("L97").Value = 1486.00002
$ws.Range("N97").Value = -2478.00002

$ws.Range("H112").Value = 2075.8064
$ws.Range("I112").Value = 2639.6
$ws.Range("J112").Value = 1967.3846
$ws.Range("K112").Value = 7918.799999999999
$ws.Range("L112").Value = 5902.1538
$ws.Range("M112").Value = -6810.799999999999
$ws.Range("N112").Value = -8118.1538

$ws.Range("H132").Value = 2799.3428
$ws.Range("I132").Value = 2625.6885
$ws.Range("K132").Value = 7877.065500000001
$ws.Range("M132").Value = -5347.065500000001

$ws.Range("H137").Value = 38463924
$ws.Range("J137").Value = 2633.0557
$ws.Range("L137").Value = 7899.1671
$ws.Range("N137").Value = -12999.1671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4028.5715
$ws.Range("I74").Value = 4181.5
$ws.Range("J74").Value = 3111
$ws.Range("K74").Value = 4181.5
$ws.Range("L74").Value = 3111
$ws.Range("M74").Value = -3307.5
$ws.Range("N74").Value = -4859

$ws.Range("H77").Value = 4028.5715
$ws.Range("I77").Value = 4181.5
$ws.Range("J77").Value = 3111
$ws.Range("K77").Value = 20907.5
$ws.Range("L77").Value = 15555
$ws.Range("M77").Value = -16539.5
$ws.Range("N77").Value = -24291

$ws.Range("H88").Value = 10418189
$ws.Range("I88").Value = 20833882
$ws.Range("J88").Value = 2496.125
$ws.Range("K88").Value = 20833882
$ws.Range("L88").Value = 2496.125
$ws.Range("M88").Value = -20833476
$ws.Range("N88").Value = -3308.125

$ws.Range("H91").Value = 10418189
$ws.Range("I91").Value = 20833882
$ws.Range("J91").Value = 2496.125
$ws.Range("K91").Value = 20833882
$ws.Range("L91").Value = 2496.125
$ws.Range("M91").Value = -20832478
$ws.Range("N91").Value = -5304.125

$ws.Range("H97").Value = 702.8946999999999
$ws.Range("I97").Value = 703.6111
$ws.Range("J97").Value = 690
$ws.Range("K97").Value = 703.6111
$ws.Range("L97").Value = 690
$ws.Range("M97").Value = -207.6111
$ws.Range("N97").Value = -1682

$ws.Range("H110").Value = 50004628
$ws.Range("I110").Value = 66670508
$ws.Range("J110").Value = 6995.4
$ws.Range("K110").Value = 66670508
$ws.Range("L110").Value = 6995.4
$ws.Range("M110").Value = -66668463
$ws.Range("N110").Value = -11085.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 13160624
$ws.Range("J86").Value = 2876.6
$ws.Range("L86").Value = 2876.6
$ws.Range("N86").Value = -5122.6

$ws.Range("H89").Value = 13160624
$ws.Range("J89").Value = 2876.6
$ws.Range("L89").Value = 14383
$ws.Range("N89").Value = -25615

$ws.Range("H94").Value = 3522.64
$ws.Range("I94").Value = 3498.55
$ws.Range("K94").Value = 3498.55
$ws.Range("M94").Value = -3047.55

$ws.Range("H134").Value = 1341.5834
$ws.Range("I134").Value = 1204.4546
$ws.Range("K134").Value = 3613.3638
$ws.Range("M134").Value = -1078.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 780.3
$ws.Range("I22").Value = 677.5294
$ws.Range("K22").Value = 677.5294
$ws.Range("M22").Value = -327.5294

$ws.Range("H31").Value = 6593.5405
$ws.Range("I31").Value = 61724.5
$ws.Range("K31").Value = 61724.5
$ws.Range("M31").Value = -61429.5

$ws.Range("H34").Value = 6593.5405
$ws.Range("I34").Value = 61724.5
$ws.Range("K34").Value = 61724.5
$ws.Range("M34").Value = -61522.5

$ws.Range("H99").Value = 2469.353
$ws.Range("J99").Value = 2748
$ws.Range("L99").Value = 2748
$ws.Range("N99").Value = -5744

$ws.Range("H105").Value = 2616.5
$ws.Range("I105").Value = 2327.5
$ws.Range("J105").Value = 3050
$ws.Range("K105").Value = 2327.5
$ws.Range("L105").Value = 3050
$ws.Range("M105").Value = -580.5
$ws.Range("N105").Value = -6544

$ws.Range("H126").Value = 2469.353
$ws.Range("J126").Value = 2748
$ws.Range("L126").Value = 8244
$ws.Range("N126").Value = -13184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11773299
$ws.Range("I4").Value = 20179860
$ws.Range("J4").Value = 4113.6333
$ws.Range("K4").Value = 60539580
$ws.Range("L4").Value = 12340.8999
$ws.Range("M4").Value = -60539468
$ws.Range("N4").Value = -12564.8999

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = ""
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = 0

$ws.Range("H129").Value = 2164.5386
$ws.Range("I129").Value = 1234.4
$ws.Range("K129").Value = 3703.2
$ws.Range("M129").Value = 1296.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3732
$ws.Range("I80").Value = 3895.5334
$ws.Range("K80").Value = 3895.5334
$ws.Range("M80").Value = -2897.5334

$ws.Range("H83").Value = 3732
$ws.Range("I83").Value = 3895.5334
$ws.Range("K83").Value = 19477.667
$ws.Range("M83").Value = -14485.667

$ws.Range("H92").Value = 9562.625
$ws.Range("J92").Value = 9562.625
$ws.Range("L92").Value = 9562.625
$ws.Range("N92").Value = -13306.625

$ws.Range("H132").Value = 6461.381
$ws.Range("I132").Value = 5440.3076
$ws.Range("J132").Value = 8120.625
$ws.Range("K132").Value = 16320.9228
$ws.Range("L132").Value = 24361.875
$ws.Range("M132").Value = -13790.9228
$ws.Range("N132").Value = -29421.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 415.92856
$ws.Range("I16").Value = 472.83334
$ws.Range("J16").Value = 74.5
$ws.Range("K16").Value = 472.83334
$ws.Range("L16").Value = 74.5
$ws.Range("M16").Value = -302.83334
$ws.Range("N16").Value = -414.5

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = ""

$ws.Range("H61").Value = 1315.909
$ws.Range("I61").Value = 1117.6
$ws.Range("K61").Value = 1117.6
$ws.Range("M61").Value = -915.5999999999999

$ws.Range("H87").Value = 90000
$ws.Range("J87").Value = 90000
$ws.Range("L87").Value = 90000
$ws.Range("N87").Value = -92246

$ws.Range("H90").Value = 90000
$ws.Range("J90").Value = 90000
$ws.Range("L90").Value = 270000
$ws.Range("N90").Value = -281232

$ws.Range("H113").Value = 1315.909
$ws.Range("I113").Value = 1117.6
$ws.Range("K113").Value = 1117.6
$ws.Range("M113").Value = 1052.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 5500
$ws.Range("J49").Value = 5500
$ws.Range("L49").Value = 5500
$ws.Range("N49").Value = -5960

$ws.Range("H81").Value = 8005580.5
$ws.Range("J81").Value = 22233772
$ws.Range("L81").Value = 44467544
$ws.Range("N81").Value = -44469666

$ws.Range("H84").Value = 8005580.5
$ws.Range("J84").Value = 22233772
$ws.Range("L84").Value = 222337720
$ws.Range("N84").Value = -222348328
